# Automatische sync: 2025-06-17 11:57:38
#
# Appends two new incoming-mail log rows to the "Logs" sheet, extends the
# conditional-formatting ranges to cover them, recomputes the "Dashboard"
# category summary (now sorted by count, descending) and widens the bar
# chart's category/value series references to match.

$wb = $excel.ActiveWorkbook

# xlPasteValues - used below to drop a multi-line string into a cell
# without Excel's "auto-fit row height for wrapped/multi-line text"
# kicking in (a plain `.Value =` assignment of a string containing
# embedded newlines triggers that autofit here, which a plain text
# entry in the source data never did).
$xlPasteValues = -4163

function Set-PlainTextValue($range, [string]$text) {
    $scratchRow = 1000
    $scratch = $range.Worksheet.Cells.Item($scratchRow, 26)   # column Z, far below real data
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial($xlPasteValues)
    $range.Worksheet.Rows.Item($scratchRow).Delete()
}

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append row 6 (Afmelding) and row 7 (Informatieaanvraag)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A6").Value = "Afmelding nieuwsbrief"
$logs.Range("B6").Value = "mailmind.test@zohomail.eu"
$logs.Range("C6").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D6").Value = "Afmelding"
$logs.Range("F6").Value = "2025-06-17 10:58:37"
$logs.Range("G6").Value = "Nee"

$logs.Range("A7").Value = "Re: Wat zijn jullie openingstijden?"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
Set-PlainTextValue $logs.Range("C7") "Beste,`nBedankt voor je interesse. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 uur tot 18:00 uur. Op zaterdag zijn wij geopend van 10:00 uur tot 16:00 uur. Voor eventuele feestdagen en afwijkende openingstijden, adviseer ik onze website te raadplegen.`nMet vriendelijke groet, [Jouw naam]"
$logs.Range("D7").Value = "Informatieaanvraag"
Set-PlainTextValue $logs.Range("E7") "Beste,`nHartelijk dank voor uw interesse. Onze openingstijden zijn van maandag tot en met vrijdag van 9:00 uur tot 18:00 uur en op zaterdag van 10:00 uur tot 16:00 uur. Voor eventuele feestdagen en afwijkende openingstijden kunt u onze website raadplegen.`nMet vriendelijke groet,`n[Jouw naam]"
$logs.Range("F7").Value = "2025-06-17 10:58:40"
$logs.Range("G7").Value = "Ja"

# Extend the conditional-formatting "applies to" ranges from rows 2-5 to 2-7
$catRules = $logs.Range("D2:D5").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D7"))
}

$answeredRules = $logs.Range("G2:G5").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G7"))
}

# ---------------------------------------------------------------------
# 2. "Dashboard" sheet: recompute category counts, sorted by count desc
#    Overig=1, Klacht=1, Bestelling=1, Informatieaanvraag=2 (was 1, now
#    +1 from new row 7), plus the brand-new "Afmelding"=1 category.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Informatieaanvraag"
$dash.Range("B2").Value = 2
$dash.Range("A3").Value = "Overig"
$dash.Range("B3").Value = 1
$dash.Range("A4").Value = "Klacht"
$dash.Range("B4").Value = 1
$dash.Range("A5").Value = "Bestelling"
$dash.Range("B5").Value = 1
$dash.Range("A6").Value = "Afmelding"
$dash.Range("B6").Value = 1

# ---------------------------------------------------------------------
# 3. Chart on "Dashboard": widen category/value series ranges from
#    row 5 to row 6 to include the new "Afmelding" bucket.
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$6,'Dashboard'!`$B`$2:`$B`$6,1)"
